$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post previously at row 672 ("「何事も永遠には続かない」") was removed.
# Delete that entire row so every following row shifts up by one
# (673 -> 672, 674 -> 673, ... 842 -> 841), matching the new used range A1:C841.
$ws.Rows(672).Delete()
